# reemplazando caldas por manizales
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename "Caldas" -> "Manizales" (row label + column header) ---
$ws.Range("A8").Value = "Manizales"
$ws.Range("H1").Value = "Manizales"

# --- Row 2 (Armenia) : new Manizales distance ---
$ws.Range("H2").Value = 96
$ws.Range("I2").Value = 268

# --- Row 3 (Barranquilla) ---
$ws.Range("H3").Value = 894
$ws.Range("I3").Value = 702

# --- Row 4 (Bogota D.C.) ---
$ws.Range("H4").Value = 302
$ws.Range("I4").Value = 443

# --- Row 5 (Bucaramanga) ---
$ws.Range("H5").Value = 508
$ws.Range("I5").Value = 404

# --- Row 6 (Cartagena) ---
$ws.Range("H6").Value = 831
$ws.Range("I6").Value = 640

# --- Row 7 (Cucuta) ---
$ws.Range("H7").Value = 694
$ws.Range("I7").Value = 598

# --- Row 8 (Manizales' own row - distances to everyone else) ---
$ws.Range("B8").Value = 96
$ws.Range("C8").Value = 894
$ws.Range("D8").Value = 302
$ws.Range("E8").Value = 508
$ws.Range("F8").Value = 831
$ws.Range("G8").Value = 694
$ws.Range("I8").Value = 194
$ws.Range("J8").Value = 600
$ws.Range("K8").Value = 258
$ws.Range("L8").Value = 54
$ws.Range("M8").Value = 424

# --- Row 9 (Medellin) - flatten to static values, H9/Manizales distance is new ---
$ws.Range("B9").Value = 268
$ws.Range("C9").Value = 702
$ws.Range("D9").Value = 443
$ws.Range("E9").Value = 404
$ws.Range("F9").Value = 640
$ws.Range("G9").Value = 598
$ws.Range("H9").Value = 194
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 405
$ws.Range("K9").Value = 425
$ws.Range("L9").Value = 215
$ws.Range("M9").Value = 802
$ws.Range("N9").Value = 702
$ws.Range("O9").Value = 425
$ws.Range("P9").Value = 745

# --- Row 10 (Monteria) ---
$ws.Range("H10").Value = 600
$ws.Range("I10").Value = 405

# --- Row 11 (Palmira) ---
$ws.Range("H11").Value = 258
$ws.Range("I11").Value = 425

# --- Row 12 (Pereira) ---
$ws.Range("H12").Value = 54
$ws.Range("I12").Value = 215

# --- Row 13 (Pasto) ---
$ws.Range("H13").Value = 424
$ws.Range("I13").Value = 802

# --- Row 14 (Soledad) : I14 loses formula, keeps old cached value ---
$ws.Range("I14").Value = 702

# --- Row 15 (Tulua) : I15 loses formula, keeps old cached value ---
$ws.Range("I15").Value = 425

# --- Row 16 (Valledupar) ---
$ws.Range("H16").Value = 860
$ws.Range("I16").Value = 745

# --- Selection moves to M20 ---
$ws.Range("M20").Select()
